$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text field corrections (Razon social / Nombre Fantasia): replace commas with periods ---
$textChanges = @(
    @('E36', 'RAMIREZ CLAUDIA. RAMIREZ CESAR Y RAMIREZ VERONICA SH'),
    @('E86', 'FERNANDEZ. MARIO HUGO'),
    @('E88', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F88', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E109', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('F109', 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'),
    @('E89', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'),
    @('E122', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'),
    @('E168', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'),
    @('E182', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'),
    @('E218', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'),
    @('E124', 'RICCOTTI. MARIANA EDITH'),
    @('F135', 'MERCANZINI. GASTON ARIEL')
)

foreach ($pair in $textChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# --- Numeric amount corrections (convert Spanish-formatted text '1.234,56' to '1234.56') ---
$numberChanges = @(
    @('H2', '7950.00'),
    @('H3', '31749.00'),
    @('H4', '53400.00'),
    @('H5', '2692750.00'),
    @('H6', '43620.50'),
    @('H7', '60396.00'),
    @('H8', '2147.00'),
    @('H9', '119991.43'),
    @('H10', '3006.00'),
    @('H11', '137.80'),
    @('H12', '50.00'),
    @('H13', '6500.00'),
    @('H14', '96869.60'),
    @('H15', '225259.63'),
    @('H16', '624.00'),
    @('H17', '33793.22'),
    @('H18', '15450.00'),
    @('H19', '961.50'),
    @('H20', '24119.22'),
    @('H21', '11972.97'),
    @('H22', '4800.00'),
    @('H23', '750.00'),
    @('H156', '750.00'),
    @('H24', '1800.00'),
    @('H214', '1800.00'),
    @('H25', '12888.00'),
    @('H26', '75.50'),
    @('H27', '54.29'),
    @('H28', '80.00'),
    @('H72', '80.00'),
    @('H74', '80.00'),
    @('H81', '80.00'),
    @('H29', '950.00'),
    @('H155', '950.00'),
    @('H30', '11745.55'),
    @('H31', '210.00'),
    @('H32', '6378.15'),
    @('H33', '128.00'),
    @('H34', '474.80'),
    @('H35', '11.00'),
    @('H36', '7596.00'),
    @('H37', '72.00'),
    @('H38', '17093.00'),
    @('H39', '674.82'),
    @('H40', '435.86'),
    @('H41', '2158.00'),
    @('H42', '27119.67'),
    @('H43', '76722.01'),
    @('H44', '48.00'),
    @('H45', '211.31'),
    @('H46', '25.00'),
    @('H47', '3240.00'),
    @('H48', '6589.99'),
    @('H49', '24.80'),
    @('H50', '563.19'),
    @('H51', '4160.00'),
    @('H52', '700.00'),
    @('H53', '22900.00'),
    @('H54', '1659.90'),
    @('H55', '2347.04'),
    @('H56', '2205.00'),
    @('H57', '3838.75'),
    @('H58', '22.00'),
    @('H59', '5535.82'),
    @('H60', '6345.24'),
    @('H61', '124.80'),
    @('H62', '3100.00'),
    @('H63', '17500.00'),
    @('H64', '1280.68'),
    @('H65', '195.00'),
    @('H66', '3522.00'),
    @('H67', '508.50'),
    @('H68', '260.00'),
    @('H69', '13604.20'),
    @('H70', '1700.00'),
    @('H71', '695.00'),
    @('H73', '104.56'),
    @('H75', '295.00'),
    @('H76', '3280.00'),
    @('H77', '8.00'),
    @('H78', '7240.00'),
    @('H79', '7360.00'),
    @('H161', '7360.00'),
    @('H80', '7890.00'),
    @('H82', '150.00'),
    @('H83', '17760.00'),
    @('H84', '60770.00'),
    @('H85', '220.00'),
    @('H86', '330.00'),
    @('H87', '15038.00'),
    @('H88', '784.00'),
    @('H89', '7568.00'),
    @('H90', '305.00'),
    @('H91', '110.00'),
    @('H92', '375.50'),
    @('H93', '55.00'),
    @('H94', '19500.00'),
    @('H95', '19000.00'),
    @('H96', '222570.00'),
    @('H97', '1730.13'),
    @('H98', '888.22'),
    @('H99', '118.37'),
    @('H100', '270.00'),
    @('H101', '435.00'),
    @('H102', '1270.69'),
    @('H103', '342.40'),
    @('H104', '3350.00'),
    @('H105', '77.23'),
    @('H106', '149.00'),
    @('H107', '43.00'),
    @('H108', '78.00'),
    @('H109', '124.00'),
    @('H110', '147.00'),
    @('H111', '552.55'),
    @('H112', '131.85'),
    @('H113', '60.00'),
    @('H114', '540.00'),
    @('H115', '2066.05'),
    @('H116', '54.00'),
    @('H117', '26.45'),
    @('H118', '6898.22'),
    @('H119', '770.00'),
    @('H120', '50.10'),
    @('H121', '490.00'),
    @('H122', '40.00'),
    @('H123', '172.86'),
    @('H124', '1000.00'),
    @('H137', '1000.00'),
    @('H147', '1000.00'),
    @('H151', '1000.00'),
    @('H125', '140.00'),
    @('H126', '1780.00'),
    @('H127', '6082.10'),
    @('H128', '2080.00'),
    @('H129', '889.00'),
    @('H130', '5213.77'),
    @('H131', '1622.00'),
    @('H132', '4608.00'),
    @('H133', '10924.00'),
    @('H134', '1060.00'),
    @('H135', '6000.00'),
    @('H136', '2000.00'),
    @('H158', '2000.00'),
    @('H138', '1411.18'),
    @('H139', '2969.00'),
    @('H140', '1749.80'),
    @('H141', '225.00'),
    @('H142', '14417.00'),
    @('H143', '131132.40'),
    @('H144', '6540.00'),
    @('H145', '1600.00'),
    @('H146', '3000.00'),
    @('H148', '20866.20'),
    @('H149', '384.00'),
    @('H150', '600.00'),
    @('H152', '5000.00'),
    @('H153', '6696.30'),
    @('H154', '1500.00'),
    @('H216', '1500.00'),
    @('H157', '3530.00'),
    @('H159', '200.00'),
    @('H160', '290.00'),
    @('H162', '1300.00'),
    @('H163', '530.00'),
    @('H164', '340.00'),
    @('H165', '968.00'),
    @('H166', '45090.00'),
    @('H167', '250.00'),
    @('H168', '2665.00'),
    @('H169', '357.51'),
    @('H170', '1775.00'),
    @('H171', '178.36'),
    @('H172', '3492.00'),
    @('H173', '800.00'),
    @('H174', '6.96'),
    @('H175', '10813.00'),
    @('H176', '390.00'),
    @('H177', '760.50'),
    @('H178', '25750.00'),
    @('H179', '815.48'),
    @('H180', '12502.00'),
    @('H181', '967.00'),
    @('H182', '21768.00'),
    @('H183', '6363.29'),
    @('H184', '47.54'),
    @('H185', '236.00'),
    @('H186', '63.00'),
    @('H187', '63.40'),
    @('H188', '25176.50'),
    @('H189', '1589.68'),
    @('H190', '4770.68'),
    @('H191', '180.00'),
    @('H192', '497.12'),
    @('H193', '799.80'),
    @('H194', '2094.24'),
    @('H195', '2883.47'),
    @('H196', '3520.00'),
    @('H197', '7200.00'),
    @('H198', '1948.25'),
    @('H199', '5830.00'),
    @('H200', '657285.86'),
    @('H201', '31435.00'),
    @('H202', '400.00'),
    @('H203', '67100.00'),
    @('H204', '216500.00'),
    @('H205', '52000.00'),
    @('H206', '36000.00'),
    @('H207', '20000.00'),
    @('H208', '50000.00'),
    @('H209', '24000.00'),
    @('H210', '223000.00'),
    @('H213', '223000.00'),
    @('H211', '7000.00'),
    @('H212', '111500.00'),
    @('H215', '564778.65'),
    @('H217', '121482.10'),
    @('H218', '2270.00'),
    @('H219', '120.00'),
    @('H220', '278449.38'),
    @('H221', '3236.30'),
    @('H222', '19200.00'),
    @('H223', '300.00'),
    @('H224', '9500.00'),
    @('H225', '118000.00'),
    @('H226', '231989.00'),
    @('H227', '212500.00'),
    @('H228', '148000.00')
)

foreach ($pair in $numberChanges) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = '@'
    $cell.Value = $pair[1]
    $cell.Style = 'Normal'
}

Write-Output 'done'